$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Version value (row 3, column B) from 0.4.0 to 0.7.0
$ws.Cells.Item(3, 2).Value = "0.7.0"

# Remove the "Jurisdiction" / "Chile" row entirely (row 11), shifting rows below up
$ws.Rows.Item(11).Delete()
